$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: Price (column D) and Volume(1h) (column E) cells are stored as text in
# this sheet. Some new prices look like plain numbers (e.g. "1.00", "583.94"),
# so a leading apostrophe (quote-prefix) is used to force Excel to keep those
# as text instead of auto-converting them to numeric values.

$ws.Range("D2").Value = '67.053.31'
$ws.Range("E2").Value = '  -1.73%  '

$ws.Range("D3").Value = '2.486.98'
$ws.Range("E3").Value = '  -1.84%  '

$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").Value = '''583.94'
$ws.Range("E5").Value = '  -1.70%  '

$ws.Range("D6").Value = '''165.44'
$ws.Range("E6").Value = '  -7.60%  '

$ws.Range("E7").Value = '  +0.17%  '

$ws.Range("D8").Value = '''0.516'
$ws.Range("E8").Value = '  -3.00%  '

$ws.Range("D9").Value = '2.490.95'
$ws.Range("E9").Value = '  -1.62%  '

$ws.Range("E10").Value = '  -4.51%  '

$ws.Range("E11").Value = '  +0.10%  '

$ws.Range("E12").Value = '  -2.05%  '

$ws.Range("D13").Value = '''4.94'
$ws.Range("E13").Value = '  -3.44%  '

$ws.Range("D14").Value = '''25.89'
$ws.Range("E14").Value = '  -3.75%  '

$ws.Range("D15").Value = '2.944.89'
$ws.Range("E15").Value = '  -1.67%  '

$ws.Range("D16").Value = '''0.0000173'
$ws.Range("E16").Value = '  -3.98%  '

$ws.Range("D17").Value = '67.033.40'
$ws.Range("E17").Value = '  -1.63%  '

$ws.Range("D18").Value = '2.472.22'
$ws.Range("E18").Value = '  -2.71%  '

$ws.Range("D19").Value = '''11.60'
$ws.Range("E19").Value = '  +0.45%  '

$ws.Range("E20").Value = '  -2.11%  '

$ws.Range("D21").Value = '''359.35'
$ws.Range("E21").Value = '  -2.27%  '

$ws.Range("D22").Value = '''4.11'
$ws.Range("E22").Value = '  -2.31%  '

$ws.Range("D23").Value = '''4.41'
$ws.Range("E23").Value = '  -6.59%  '

$ws.Range("E24").Value = '  +0.10%  '

$ws.Range("E25").Value = '  -0.36%  '

$ws.Range("D26").Value = '''1.84'
$ws.Range("E26").Value = '  -4.99%  '

$ws.Range("E27").Value = '  -7.93%  '

$ws.Range("D28").Value = '''0.991'
$ws.Range("E28").Value = '  -0.68%  '

$ws.Range("D29").Value = '2.614.31'

$ws.Range("D30").Value = '0.0₃0933'
$ws.Range("E30").Value = '  -6.72%  '

$ws.Range("D31").Value = '''8.04'
$ws.Range("E31").Value = '  -3.23%  '

$ws.Range("D32").Value = '''498.14'
$ws.Range("E32").Value = '  -8.04%  '

$ws.Range("E33").Value = '  -2.41%  '

$ws.Range("E34").Value = '  -5.60%  '

$ws.Range("D35").Value = '''1.00'
$ws.Range("E35").Value = '  +0.27%  '

$ws.Range("E36").Value = '  -2.46%  '

$ws.Range("E37").Value = '  +1.50%  '

$ws.Range("D38").Value = '''19.01'
$ws.Range("E38").Value = '  +0.67%  '

$ws.Range("E39").Value = '  -3.64%  '

$ws.Range("D40").Value = '''18.57'
$ws.Range("E40").Value = '  -0.65%  '

$ws.Range("E41").Value = '  -4.88%  '

$ws.Range("D42").Value = '''4.93'
$ws.Range("E42").Value = '  -5.61%  '

$ws.Range("E43").Value = '  -5.89%  '

$ws.Range("E44").Value = '  +0.05%  '

$ws.Range("D45").Value = '''2.46'
$ws.Range("E45").Value = '  -4.70%  '

$ws.Range("D46").Value = '''39.33'
$ws.Range("E46").Value = '  -1.49%  '

$ws.Range("D47").Value = '''141.82'
$ws.Range("E47").Value = '  -3.95%  '

$ws.Range("E48").Value = '  -3.00%  '

$ws.Range("D49").Value = '''0.537'
$ws.Range("E49").Value = '  -4.57%  '

$ws.Range("D50").Value = '0.0₆0263'
$ws.Range("E50").Value = '  -5.73%  '

$ws.Range("E51").Value = '  -3.70%  '
